# Update the date line and all multiplication problems in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2023-10-10 Tuesday", "2023-10-11 Wednesday"),
    @("27×79=", "58×89="),
    @("77×53=", "57×66="),
    @("61×93=", "44×36="),
    @("62×66=", "89×55="),
    @("23×65=", "50×77="),
    @("26×33=", "79×18="),
    @("63×60=", "53×48="),
    @("15×41=", "80×51="),
    @("71×18=", "42×15="),
    @("13×74=", "60×22="),
    @("69×83=", "72×50="),
    @("92×34=", "84×11="),
    @("89×41=", "53×55="),
    @("98×48=", "98×90="),
    @("42×67=", "93×40="),
    @("12×60=", "19×67="),
    @("83×97=", "92×13="),
    @("87×47=", "84×56="),
    @("28×80=", "24×49="),
    @("77×22=", "69×69="),
    @("76×66=", "13×72="),
    @("24×30=", "63×41="),
    @("72×41=", "71×55="),
    @("67×20=", "97×40="),
    @("13×61=", "66×27=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
